$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_8_9_24"
$ws.Cells.Item(2, 2).Value = 0.0142660343961194
$ws.Cells.Item(2, 3).Value = -0.686191144538788
$ws.Cells.Item(2, 4).Value = -4.214568390048567
$ws.Cells.Item(2, 5).Value = -2.276669917300878
$ws.Cells.Item(2, 6).Value = 1.090916991233826
$ws.Cells.Item(2, 7).Value = 1.577954411506653
$ws.Cells.Item(2, 8).Value = 5.861174583435059
$ws.Cells.Item(2, 9).Value = 3.593586444854736
$ws.Cells.Item(3, 1).Value = "model_8_9_23"
$ws.Cells.Item(3, 2).Value = 0.01428869368752717
$ws.Cells.Item(3, 3).Value = -0.6861650363286458
$ws.Cells.Item(3, 4).Value = -4.214434666534983
$ws.Cells.Item(3, 5).Value = -2.276594570291274
$ws.Cells.Item(3, 6).Value = 1.09089195728302
$ws.Cells.Item(3, 7).Value = 1.577929973602295
$ws.Cells.Item(3, 8).Value = 5.861023902893066
$ws.Cells.Item(3, 9).Value = 3.593503952026367
$ws.Cells.Item(4, 1).Value = "model_8_9_22"
$ws.Cells.Item(4, 2).Value = 0.01453367816956652
$ws.Cells.Item(4, 3).Value = -0.6858589898143499
$ws.Cells.Item(4, 4).Value = -4.213031990655103
$ws.Cells.Item(4, 5).Value = -2.275779808598656
$ws.Cells.Item(4, 6).Value = 1.090620756149292
$ws.Cells.Item(4, 7).Value = 1.577643752098083
$ws.Cells.Item(4, 8).Value = 5.859447479248047
$ws.Cells.Item(4, 9).Value = 3.592610597610474
$ws.Cells.Item(5, 1).Value = "model_8_9_21"
$ws.Cells.Item(5, 2).Value = 0.01503948113420961
$ws.Cells.Item(5, 3).Value = -0.6857454979163053
$ws.Cells.Item(5, 4).Value = -4.209651986313005
$ws.Cells.Item(5, 5).Value = -2.274097480188
$ws.Cells.Item(5, 6).Value = 1.090061068534851
$ws.Cells.Item(5, 7).Value = 1.577537417411804
$ws.Cells.Item(5, 8).Value = 5.855648517608643
$ws.Cells.Item(5, 9).Value = 3.590765237808228
$ws.Cells.Item(6, 1).Value = "model_8_9_20"
$ws.Cells.Item(6, 2).Value = 0.01532256653503938
$ws.Cells.Item(6, 3).Value = -0.6850832124894193
$ws.Cells.Item(6, 4).Value = -4.208315627813604
$ws.Cells.Item(6, 5).Value = -2.273156267033043
$ws.Cells.Item(6, 6).Value = 1.089747786521912
$ws.Cells.Item(6, 7).Value = 1.57691764831543
$ws.Cells.Item(6, 8).Value = 5.854146480560303
$ws.Cells.Item(6, 9).Value = 3.589732646942139
$ws.Cells.Item(7, 1).Value = "model_8_9_19"
$ws.Cells.Item(7, 2).Value = 0.01536985982920469
$ws.Cells.Item(7, 3).Value = -0.6856196381725685
$ws.Cells.Item(7, 4).Value = -4.207487477743493
$ws.Cells.Item(7, 5).Value = -2.272998948391852
$ws.Cells.Item(7, 6).Value = 1.089695334434509
$ws.Cells.Item(7, 7).Value = 1.577419757843018
$ws.Cells.Item(7, 8).Value = 5.853215217590332
$ws.Cells.Item(7, 9).Value = 3.589560270309448
$ws.Cells.Item(8, 1).Value = "model_8_9_18"
$ws.Cells.Item(8, 2).Value = 0.01719689771016719
$ws.Cells.Item(8, 3).Value = -0.6834152047994186
$ws.Cells.Item(8, 4).Value = -4.19695563060304
$ws.Cells.Item(8, 5).Value = -2.266920702312053
$ws.Cells.Item(8, 6).Value = 1.087673425674438
$ws.Cells.Item(8, 7).Value = 1.575356602668762
$ws.Cells.Item(8, 8).Value = 5.841377258300781
$ws.Cells.Item(8, 9).Value = 3.582894325256348
$ws.Cells.Item(9, 1).Value = "model_8_9_17"
$ws.Cells.Item(9, 2).Value = 0.02067111599315974
$ws.Cells.Item(9, 3).Value = -0.6745531689366748
$ws.Cells.Item(9, 4).Value = -4.181296666737246
$ws.Cells.Item(9, 5).Value = -2.255367047955548
$ws.Cells.Item(9, 6).Value = 1.083828449249268
$ws.Cells.Item(9, 7).Value = 1.567063570022583
$ws.Cells.Item(9, 8).Value = 5.823777198791504
$ws.Cells.Item(9, 9).Value = 3.570223093032837
$ws.Cells.Item(10, 1).Value = "model_8_9_16"
$ws.Cells.Item(10, 2).Value = 0.02872167010350268
$ws.Cells.Item(10, 3).Value = -0.6674379704919049
$ws.Cells.Item(10, 4).Value = -4.132476387278895
$ws.Cells.Item(10, 5).Value = -2.228605430493992
$ws.Cells.Item(10, 6).Value = 1.074918985366821
$ws.Cells.Item(10, 7).Value = 1.560405015945435
$ws.Cells.Item(10, 8).Value = 5.768902778625488
$ws.Cells.Item(10, 9).Value = 3.540873289108276
$ws.Cells.Item(11, 1).Value = "model_8_9_15"
$ws.Cells.Item(11, 2).Value = 0.02902105315443504
$ws.Cells.Item(11, 3).Value = -0.6702790608517275
$ws.Cells.Item(11, 4).Value = -4.127738488844255
$ws.Cells.Item(11, 5).Value = -2.22760497368994
$ws.Cells.Item(11, 6).Value = 1.07458758354187
$ws.Cells.Item(11, 7).Value = 1.563063740730286
$ws.Cells.Item(11, 8).Value = 5.763577461242676
$ws.Cells.Item(11, 9).Value = 3.539776086807251
$ws.Cells.Item(12, 1).Value = "model_8_9_14"
$ws.Cells.Item(12, 2).Value = 0.03196150667886777
$ws.Cells.Item(12, 3).Value = -0.6699625434963952
$ws.Cells.Item(12, 4).Value = -4.107700879694351
$ws.Cells.Item(12, 5).Value = -2.217797871397296
$ws.Cells.Item(12, 6).Value = 1.071333289146423
$ws.Cells.Item(12, 7).Value = 1.562767624855042
$ws.Cells.Item(12, 8).Value = 5.741055488586426
$ws.Cells.Item(12, 9).Value = 3.529020309448242
$ws.Cells.Item(13, 1).Value = "model_8_9_13"
$ws.Cells.Item(13, 2).Value = 0.04617166806055084
$ws.Cells.Item(13, 3).Value = -0.6495430460857872
$ws.Cells.Item(13, 4).Value = -4.028776778750252
$ws.Cells.Item(13, 5).Value = -2.170509601667613
$ws.Cells.Item(13, 6).Value = 1.055606961250305
$ws.Cells.Item(13, 7).Value = 1.543658852577209
$ws.Cells.Item(13, 8).Value = 5.652344703674316
$ws.Cells.Item(13, 9).Value = 3.477158308029175
$ws.Cells.Item(14, 1).Value = "model_8_9_12"
$ws.Cells.Item(14, 2).Value = 0.06314268413903867
$ws.Cells.Item(14, 3).Value = -0.6284488187897554
$ws.Cells.Item(14, 4).Value = -3.930722023618835
$ws.Cells.Item(14, 5).Value = -2.113689176375491
$ws.Cells.Item(14, 6).Value = 1.036824941635132
$ws.Cells.Item(14, 7).Value = 1.523918628692627
$ws.Cells.Item(14, 8).Value = 5.542130947113037
$ws.Cells.Item(14, 9).Value = 3.41484260559082
$ws.Cells.Item(15, 1).Value = "model_8_9_11"
$ws.Cells.Item(15, 2).Value = 0.09176154994266783
$ws.Cells.Item(15, 3).Value = -0.5751264689790598
$ws.Cells.Item(15, 4).Value = -3.78086018264193
$ws.Cells.Item(15, 5).Value = -2.01732533598787
$ws.Cells.Item(15, 6).Value = 1.005152344703674
$ws.Cells.Item(15, 7).Value = 1.474019050598145
$ws.Cells.Item(15, 8).Value = 5.373686790466309
$ws.Cells.Item(15, 9).Value = 3.309158325195312
$ws.Cells.Item(16, 1).Value = "model_8_9_10"
$ws.Cells.Item(16, 2).Value = 0.1355222151019148
$ws.Cells.Item(16, 3).Value = -0.5216460071840614
$ws.Cells.Item(16, 4).Value = -3.522897229432981
$ws.Cells.Item(16, 5).Value = -1.868750388823841
$ws.Cells.Item(16, 6).Value = 0.9567221403121948
$ws.Cells.Item(16, 7).Value = 1.42397141456604
$ws.Cells.Item(16, 8).Value = 5.083736419677734
$ws.Cells.Item(16, 9).Value = 3.146213293075562
$ws.Cells.Item(17, 1).Value = "model_8_9_9"
$ws.Cells.Item(17, 2).Value = 0.1688271991608291
$ws.Cells.Item(17, 3).Value = -0.4420917449086732
$ws.Cells.Item(17, 4).Value = -3.353182269017908
$ws.Cells.Item(17, 5).Value = -1.75096118199786
$ws.Cells.Item(17, 6).Value = 0.9198633432388306
$ws.Cells.Item(17, 7).Value = 1.349523782730103
$ws.Cells.Item(17, 8).Value = 4.892976760864258
$ws.Cells.Item(17, 9).Value = 3.017031669616699
$ws.Cells.Item(18, 1).Value = "model_8_9_8"
$ws.Cells.Item(18, 2).Value = 0.3264597139418017
$ws.Cells.Item(18, 3).Value = -0.04577927267177651
$ws.Cells.Item(18, 4).Value = -2.618064990786939
$ws.Cells.Item(18, 5).Value = -1.217388531265002
$ws.Cells.Item(18, 6).Value = 0.7454105615615845
$ws.Cells.Item(18, 7).Value = 0.9786506295204163
$ws.Cells.Item(18, 8).Value = 4.066704750061035
$ws.Cells.Item(18, 9).Value = 2.431852340698242
$ws.Cells.Item(19, 1).Value = "model_8_9_2"
$ws.Cells.Item(19, 2).Value = 0.4658927792872345
$ws.Cells.Item(19, 3).Value = 0.2158803854747712
$ws.Cells.Item(19, 4).Value = -0.8992038063547936
$ws.Cells.Item(19, 5).Value = -0.2701915609206689
$ws.Cells.Item(19, 6).Value = 0.5910993218421936
$ws.Cells.Item(19, 7).Value = 0.733786940574646
$ws.Cells.Item(19, 8).Value = 2.134705066680908
$ws.Cells.Item(19, 9).Value = 1.393043279647827
$ws.Cells.Item(20, 1).Value = "model_8_9_0"
$ws.Cells.Item(20, 2).Value = 0.4880375525350504
$ws.Cells.Item(20, 3).Value = 0.7907230832902306
$ws.Cells.Item(20, 4).Value = 0.8494964772220306
$ws.Cells.Item(20, 5).Value = 0.8328746350337507
$ws.Cells.Item(20, 6).Value = 0.5665915608406067
$ws.Cells.Item(20, 7).Value = 0.1958434134721756
$ws.Cells.Item(20, 8).Value = 0.1691659390926361
$ws.Cells.Item(20, 9).Value = 0.1832896023988724
$ws.Cells.Item(21, 1).Value = "model_8_9_1"
$ws.Cells.Item(21, 2).Value = 0.4979377091613004
$ws.Cells.Item(21, 3).Value = 0.7143505074808063
$ws.Cells.Item(21, 4).Value = 0.4325449552530423
$ws.Cells.Item(21, 5).Value = 0.5972798963679282
$ws.Cells.Item(21, 6).Value = 0.5556349754333496
$ws.Cells.Item(21, 7).Value = 0.2673136293888092
$ws.Cells.Item(21, 8).Value = 0.6378194093704224
$ws.Cells.Item(21, 9).Value = 0.441670835018158
$ws.Cells.Item(22, 1).Value = "model_8_9_5"
$ws.Cells.Item(22, 2).Value = 0.5541279729083624
$ws.Cells.Item(22, 3).Value = 0.5359914217677115
$ws.Cells.Item(22, 4).Value = -1.342257762238201
$ws.Cells.Item(22, 5).Value = -0.3392685223879568
$ws.Cells.Item(22, 6).Value = 0.493448942899704
$ws.Cells.Item(22, 7).Value = 0.434223860502243
$ws.Cells.Item(22, 8).Value = 2.632697582244873
$ws.Cells.Item(22, 9).Value = 1.468801379203796
$ws.Cells.Item(23, 1).Value = "model_8_9_7"
$ws.Cells.Item(23, 2).Value = 0.5602716756795878
$ws.Cells.Item(23, 3).Value = 0.1732074302757919
$ws.Cells.Item(23, 4).Value = -1.15511255981434
$ws.Cells.Item(23, 5).Value = -0.412891454298419
$ws.Cells.Item(23, 6).Value = 0.4866497218608856
$ws.Cells.Item(23, 7).Value = 0.7737206816673279
$ws.Cells.Item(23, 8).Value = 2.422346115112305
$ws.Cells.Item(23, 9).Value = 1.549545049667358
$ws.Cells.Item(24, 1).Value = "model_8_9_4"
$ws.Cells.Item(24, 2).Value = 0.568232723531567
$ws.Cells.Item(24, 3).Value = 0.5878555574601406
$ws.Cells.Item(24, 4).Value = -1.278055649026113
$ws.Cells.Item(24, 5).Value = -0.2848730734521894
$ws.Cells.Item(24, 6).Value = 0.4778391420841217
$ws.Cells.Item(24, 7).Value = 0.3856889009475708
$ws.Cells.Item(24, 8).Value = 2.560534477233887
$ws.Cells.Item(24, 9).Value = 1.409144759178162
$ws.Cells.Item(25, 1).Value = "model_8_9_3"
$ws.Cells.Item(25, 2).Value = 0.5842529227831739
$ws.Cells.Item(25, 3).Value = 0.6581374640798341
$ws.Cells.Item(25, 4).Value = -1.212292564446782
$ws.Cells.Item(25, 5).Value = -0.2214084287548026
$ws.Cells.Item(25, 6).Value = 0.4601095020771027
$ws.Cells.Item(25, 7).Value = 0.3199183642864227
$ws.Cells.Item(25, 8).Value = 2.486616611480713
$ws.Cells.Item(25, 9).Value = 1.339541912078857
$ws.Cells.Item(26, 1).Value = "model_8_9_6"
$ws.Cells.Item(26, 2).Value = 0.7408732073321167
$ws.Cells.Item(26, 3).Value = 0.3049921024759104
$ws.Cells.Item(26, 4).Value = 0.04687200154060278
$ws.Cells.Item(26, 5).Value = 0.226351006317564
$ws.Cells.Item(26, 6).Value = 0.2867770195007324
$ws.Cells.Item(26, 7).Value = 0.6503953337669373
$ws.Cells.Item(26, 8).Value = 1.071315765380859
$ws.Cells.Item(26, 9).Value = 0.8484756350517273
